$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: add P1 = 14, Q1 = 15, using same style (bold/centered/bordered) as O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For data rows 2-25: swap I<->K and M<->O values, then add P=2, Q=2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # I
    $kVal = $ws.Cells.Item($r, 11).Value2  # K
    $mVal = $ws.Cells.Item($r, 13).Value2  # M
    $oVal = $ws.Cells.Item($r, 15).Value2  # O

    $ws.Cells.Item($r, 9).Value = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value = 2      # P
    $ws.Cells.Item($r, 17).Value = 2      # Q
}
